# Rename the SN1/SN2 fatty-acid labels to FA1/FA2 (TG [M+Na]+ / [M+H]+ fix)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FA1_[FA-H]-"
$ws.Range("A3").Value = "FA2_[FA-H]-"
$ws.Range("A4").Value = "[LPL(FA1)-H]-"
$ws.Range("A5").Value = "[LPL(FA2)-H]-"
$ws.Range("A6").Value = "[LPL(FA1)-H2O-H]-"
$ws.Range("A7").Value = "[LPL(FA2)-H2O-H]-"

# Move/save the active selection to A7, matching the last-edited cell
$null = $ws.Range("A7").Select()
